# Update cryptos list (price/volume) as scraped on Sat Sep 16 11:30:58 UTC 2023
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = @{ D = "26.644.95";  E = $null }
    3  = @{ D = "1.643.08";   E = "  +0.66%  " }
    4  = @{ D = $null;        E = "  +0.01%  " }
    5  = @{ D = "215.37";     E = "  +0.89%  " }
    6  = @{ D = $null;        E = "  +1.33%  " }
    7  = @{ D = $null;        E = "  +0.08%  " }
    8  = @{ D = $null;        E = "  +0.25%  " }
    9  = @{ D = $null;        E = "  +0.74%  " }
    10 = @{ D = "19.27";      E = "  +0.52%  " }
    11 = @{ D = $null;        E = "  -0.26%  " }
    12 = @{ D = "1.871.67";   E = "  +0.61%  " }
    13 = @{ D = "1.649.74";   E = "  +1.47%  " }
    14 = @{ D = $null;        E = "  +2.05%  " }
    15 = @{ D = $null;        E = "  +1.56%  " }
    16 = @{ D = $null;        E = "  +2.94%  " }
    17 = @{ D = "26.703.52";  E = "  +0.22%  " }
    18 = @{ D = "0.0₃0743";   E = "  +0.32%  " }
    19 = @{ D = "217.14";     E = "  +0.94%  " }
    20 = @{ D = $null;        E = "  -0.06%  " }
    21 = @{ D = $null;        E = "  +0.86%  " }
    22 = @{ D = $null;        E = "  +2.34%  " }
    23 = @{ D = "9.49";       E = "  +1.50%  " }
    24 = @{ D = $null;        E = "  +13.78%  " }
    25 = @{ D = "145.68";     E = "  -1.34%  " }
    26 = @{ D = $null;        E = "  +0.13%  " }
    27 = @{ D = $null;        E = "  -0.66%  " }
    28 = @{ D = "7.19";       E = "  +4.91%  " }
    29 = @{ D = $null;        E = "  +1.42%  " }
    30 = @{ D = $null;        E = "  +2.67%  " }
    31 = @{ D = $null;        E = "  +0.80%  " }
    32 = @{ D = $null;        E = "  +2.45%  " }
    33 = @{ D = $null;        E = "  +2.66%  " }
    34 = @{ D = "1.278.14";   E = "  +4.38%  " }
    35 = @{ D = $null;        E = "  +2.86%  " }
    36 = @{ D = $null;        E = "  +4.99%  " }
    37 = @{ D = $null;        E = "  +0.44%  " }
    38 = @{ D = "0.531";      E = "  +6.05%  " }
    39 = @{ D = $null;        E = "  +3.08%  " }
    40 = @{ D = $null;        E = "  +0.04%  " }
    41 = @{ D = "0.817";      E = "  +2.53%  " }
    42 = @{ D = "2.25";       E = "  -1.68%  " }
    43 = @{ D = $null;        E = "  +2.27%  " }
    44 = @{ D = "1.782.51";   E = "  +0.65%  " }
    45 = @{ D = "91.83";      E = "  -1.21%  " }
    46 = @{ D = "59.83";      E = "  +8.44%  " }
    47 = @{ D = $null;        E = "  +1.85%  " }
    48 = @{ D = $null;        E = "  +0.66%  " }
    49 = @{ D = "7.80";       E = "  +2.34%  " }
    50 = @{ D = "0.0971";     E = "  +3.24%  " }
    51 = @{ D = "0.406";      E = "  -0.89%  " }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    if ($null -ne $vals.D) {
        # Force text format so numeric-looking strings (e.g. "215.37") aren't
        # silently coerced into numbers by Excel, matching the inlineStr/text
        # representation used in the source workbook.
        $ws.Cells.Item($row, 4).NumberFormat = "@"
        $ws.Cells.Item($row, 4).Value = $vals.D
    }
    if ($null -ne $vals.E) {
        $ws.Cells.Item($row, 5).NumberFormat = "@"
        $ws.Cells.Item($row, 5).Value = $vals.E
    }
}
